# Updated cryptos list on Mon Apr 15 21:13:17 UTC 2024 with GitHub Actions
# Applies price/volume/coin updates per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.348.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.095.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.090.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.159"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000218"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.581.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.371.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.088.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "504.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.706"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "59.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "522.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0412"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.45%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.124"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0791"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.053.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.253"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.83%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.53%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0497"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +68.84%  "
